$wb = $excel.ActiveWorkbook

# Update the "config" sheet: B4 (name -> value of "nlp/test_data_02") becomes "[tmp]/nlp/test_data_02"
# and B5 (train.iteration) becomes 28.
$configSheet = $wb.Worksheets.Item("config")
$configSheet.Range("B4").Value = "[tmp]/nlp/test_data_02"
$configSheet.Range("B5").Value = 28

# Make "config" the active sheet (it becomes the selected/active tab).
$configSheet.Activate()
